$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $findText"
    }
    return $ok
}

# --- Paragraph 42: "Git is another example..." -> "The biggest change 3rd generation..." ---
$old42 = "Git is another example of a distributed VCS, which isn’t exactly unique anymore, but has in-fact become an industry standard. After all, it’s not exactly convenient to share a machine with someone in this day and age, haha. What Git has to uniquely offer instead, is a series of powerful tools for that work just a bit differently from most other VCS. While the most important function any VCS offers can be argued to be it’s .diff function, Git actually works by taking a snapshot of the directory, creating a reference to those files in their exact form at that moment. Then when you make another commit, Git takes another snapshot and creates another reference. While this seems inefficient, Git handles these snapshots in a very clever way, checking whether a file has changed versions when you commit. When a file does not change through a commit, Git saves a link to the previous version (with no difference) instead. Unlike most other VCS that work by tracking diffs, this distinction with snapshots allows Git to offer one of it’s most powerful tools, branching."
$new42 = "The biggest change 3rd generation VCS brought about the ability to maintain independent reposts (or forks). As of today, Git is the most popular VCS in the world, and it follows the industry standard of being a distributed VCS. Git offers users a series of powerful tools, that work just a bit differently from most other VCS. While the most important function any VCS offers can be argued to be its .diff function, Git subverts that expectation by taking a snapshot of a directory, creating a reference to all files in their exact form at that moment. Then when you make another commit, Git takes another snapshot and creates another reference. From these two snapshots, the .diff is generated. While this seems inefficient when you consider most other VCS compare the files directly to get their .diff, Git handles these snapshots in a very clever way, checking whether a file has changed versions when you commit. When a file does not change through a commit, Git saves a link to the previous version (with no difference) instead. This distinction allows Git to offer one of its most powerful tools, branching."
Replace-Text $old42 $new42

# Make the "rd" in "3rd generation" (of "The biggest change 3rd generation") superscript.
$r = $d.Content
$found = $r.Find.Execute("change 3rd generation VCS brought", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rdStart = $r.Start + 8
    $rdEnd = $r.Start + 10
    $sub = $d.Range($rdStart, $rdEnd)
    if ($sub.Text -eq "rd") {
        $sub.Font.Superscript = $true
    } else {
        Write-Output "MISMATCH sub.Text=[$($sub.Text)]"
    }
} else {
    Write-Output "NOT FOUND: change 3rd generation VCS brought"
}

# --- Paragraph 44: Branching paragraph rewrite ---
$old44 = "Branching, or parallel development, allows code to be tested in isolation, preventing conflicts and user mishaps on the most sacred source code. While most modern VCS allow you to revert mistakes, Git seeks to prevent these mistakes before they ever happen! To further take control over your project, Git takes history rewriting a step further than most, allowing you to change the order of commits, change commit messages, combine or split up multiple commits, or even remove them completely! "
$new44 = "Branching, or parallel development, allows code to be tested in isolation, preventing conflicts and user mishaps on primary source code. Git’s branching is incredibly unique when compared to other VCS in that the snapshots it takes are pointers to commits, instead of stored changes. This allows the Git branching feature to be incredibly fast, and lightweight compared to other VCS which often require you to create a new copy of your source code when branching, which can be costly for large projects. To further take control over your project, Git takes history rewriting a step further than most, allowing you to change the order of commits, change commit messages, combine or split up multiple commits, or even remove them completely! "
Replace-Text $old44 $new44

# --- Paragraph 46: "Git is fairly unique..." rewrite, split around existing superscript "nd" run ---
$old46a = "Git is fairly unique in that it was an early distributed VCS, and it found widespread popularity, especially through paired systems like Github or Gitlab. Because Git is used locally, it is incredibly fast after a developer clones a repo to their local machine, allowing for easy offline work, and a secondary backup system as well. All of this seems like a lot to be offered, and it really is, but sometimes users need, or even prefer, something a little different. That’s where the 2"
$new46a = "Git is unique in that it was one of the earliest distributed VCS, and it almost immediately found widespread popularity, especially through paired systems like Github or Gitlab. Because Git is used locally, it is incredibly fast after a developer clones a repo to their local machine, allowing for easy offline work, and a secondary backup system as well. All of this seems like a lot to be offered, and it really is, but sometimes users need, or even prefer, something a little different. That is where the 2"
Replace-Text $old46a $new46a

$old46b = " most popular VCS of the modern generation comes in."
$new46b = " most popular VCS of the 3rd generation comes in."
Replace-Text $old46b $new46b

# --- Paragraph 54: "short demo" -> "demo" ---
Replace-Text "a short demo on the history of VCS!" "a demo on the history of VCS!"

Write-Output "done"
